$d = $word.ActiveDocument

$pairs = @(
    @("48÷5=", "66÷9="),
    @("96÷9=", "55÷5="),
    @("49÷5=", "22÷6="),
    @("52÷2=", "53÷6="),
    @("29÷9=", "41÷8="),
    @("27÷9=", "82÷3="),
    @("77÷8=", "39÷6="),
    @("70÷3=", "18÷5="),
    @("49÷3=", "80÷8="),
    @("78÷8=", "67÷5="),
    @("16÷5=", "78÷9="),
    @("25÷2=", "10÷2="),
    @("38÷2=", "36÷6="),
    @("77÷9=", "23÷4="),
    @("69÷8=", "63÷3="),
    @("24÷2=", "68÷6="),
    @("41÷4=", "63÷4="),
    @("11÷4=", "62÷8="),
    @("14÷9=", "11÷9="),
    @("92÷4=", "51÷4="),
    @("77÷5=", "82÷2="),
    @("28÷5=", "90÷7="),
    @("17÷6=", "12÷6="),
    @("66÷5=", "90÷9="),
    @("43÷5=", "54÷3=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
